# Append two new submission records to the "Records" sheet.
# Each element is one row in column order: Time, Amount, Purpose, Invoice, User, Status, Reviewer
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2026/1/18 02:15:53", 2026.1, "Hello World", '["/uploads/Rylan/invoice-1768673630289-649280843.webp"]', "Rylan", "REJECTED", "ADMIN"),
    @("2026/1/18 02:29:33", 1, "二次上传test", '["/uploads/Rylan/invoice-1768674560472-457905327.webp"]', "Rylan", "REJECTED", "ADMIN")
)

$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $newRows[$i]
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
}
